$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 78).
for ($row = 2; $row -le 78; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
